$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above the current row 22, shifting existing rows
# 22-32 down to 23-33 (and extending the sheet dimension to R33).
$ws.Rows(22).Insert()

# Populate the newly inserted row 22 with the new weekly price record.
$ws.Range("A22").Value = 5
$ws.Range("B22").Value = "Macroferia Regional de Talca"
$ws.Range("C22").Value = "Maule"
$ws.Range("D22").Value2 = 44488
$ws.Range("E22").Value = 7
$ws.Range("F22").Value = 300000000
$ws.Range("G22").Value = "Espárragos"
$ws.Range("H22").Value = "Verde"
$ws.Range("I22").Value = "Primera"
$ws.Range("J22").Value = 6000
$ws.Range("K22").Value = 850
$ws.Range("L22").Value = 900
$ws.Range("M22").Value = 875
$ws.Range("N22").Value = "$/kilo"
$ws.Range("O22").Value = "Provincia de Linares"
$ws.Range("P22").Value = 875
$ws.Range("Q22").Value = 1
$ws.Range("R22").Value = "Hortaliza"
